$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("A4").Value = "vergeschlossen"
$ws.Range("B4").Value = "closed, locked"

$ws.Range("A5").Value = "Freigabe"
$ws.Range("B5").Value = "release"

$ws.Range("A6").Value = "vorgange"
$ws.Range("B6").Value = "process"

$ws.Range("A7").Value = "erforderlich"
$ws.Range("B7").Value = "necessary"

$ws.Range("C7").Select()
